# Updated cryptos list values (price + 1h volume/change columns, and a coin-order
# swap for rows 9/10) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.709.48'
$ws.Range('E2').Value = '  -0.41%  '

# Row 3
$ws.Range('D3').Value = '1.869.37'
$ws.Range('E3').Value = '  -0.96%  '

# Row 4
$ws.Range('D4').Value = "'1.011"
$ws.Range('E4').Value = '  +0.36%  '

# Row 5
$ws.Range('D5').Value = "'336.37"
$ws.Range('E5').Value = '  +0.68%  '

# Row 6
$ws.Range('E6').Value = '  +0.31%  '

# Row 7
$ws.Range('D7').Value = "'0.4674"
$ws.Range('E7').Value = '  -0.85%  '

# Row 8
$ws.Range('D8').Value = "'0.3930"
$ws.Range('E8').Value = '  -0.11%  '

# Row 9
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = "'45.50"
$ws.Range('E9').Value = '  -4.48%  '

# Row 10
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.08007"
$ws.Range('E10').Value = '  -0.97%  '

# Row 11
$ws.Range('D11').Value = "'1.005"
$ws.Range('E11').Value = '  -2.25%  '

# Row 12
$ws.Range('D12').Value = "'21.91"
$ws.Range('E12').Value = '  -1.44%  '

# Row 13
$ws.Range('D13').Value = '1.880.46'
$ws.Range('E13').Value = '  -0.15%  '

# Row 14
$ws.Range('D14').Value = "'5.999"
$ws.Range('E14').Value = '  +0.24%  '

# Row 15
$ws.Range('D15').Value = "'7.269"
$ws.Range('E15').Value = '  +1.87%  '

# Row 16
$ws.Range('E16').Value = '  +0.26%  '

# Row 17
$ws.Range('D17').Value = "'88.75"
$ws.Range('E17').Value = '  +1.62%  '

# Row 18
$ws.Range('D18').Value = "'0.06740"
$ws.Range('E18').Value = '  -0.46%  '

# Row 19
$ws.Range('D19').Value = "'0.00001045"
$ws.Range('E19').Value = '  -0.44%  '

# Row 20
$ws.Range('E20').Value = '  -0.55%  '

# Row 21
$ws.Range('D21').Value = "'1.010"
$ws.Range('E21').Value = '  +0.33%  '

# Row 22
$ws.Range('D22').Value = '27.728.53'
$ws.Range('E22').Value = '  -0.41%  '

# Row 23
$ws.Range('D23').Value = "'5.489"
$ws.Range('E23').Value = '  -0.89%  '

# Row 24
$ws.Range('E24').Value = '  -0.79%  '

# Row 25
$ws.Range('D25').Value = "'2.313"
$ws.Range('E25').Value = '  -0.93%  '

# Row 26
$ws.Range('D26').Value = '2.101.25'
$ws.Range('E26').Value = '  -0.17%  '

# Row 27
$ws.Range('D27').Value = "'159.69"
$ws.Range('E27').Value = '  +0.36%  '

# Row 28
$ws.Range('D28').Value = "'19.79"
$ws.Range('E28').Value = '  -2.01%  '

# Row 29
$ws.Range('D29').Value = "'2.142"
$ws.Range('E29').Value = '  +1.54%  '

# Row 30
$ws.Range('D30').Value = "'5.458"
$ws.Range('E30').Value = '  -2.24%  '

# Row 31
$ws.Range('D31').Value = "'121.87"
$ws.Range('E31').Value = '  -0.15%  '

# Row 32
$ws.Range('D32').Value = "'0.9818"
$ws.Range('E32').Value = '  -0.18%  '

# Row 33
$ws.Range('D33').Value = "'0.09532"
$ws.Range('E33').Value = '  +0.43%  '

# Row 34
$ws.Range('D34').Value = "'3.640"
$ws.Range('E34').Value = '  +0.57%  '

# Row 35
$ws.Range('D35').Value = "'5.337"
$ws.Range('E35').Value = '  -0.57%  '

# Row 36
$ws.Range('D36').Value = "'1.338"
$ws.Range('E36').Value = '  -7.87%  '

# Row 37
$ws.Range('D37').Value = "'0.06060"
$ws.Range('E37').Value = '  -1.81%  '

# Row 38
$ws.Range('D38').Value = "'0.02237"
$ws.Range('E38').Value = '  -1.43%  '

# Row 39
$ws.Range('D39').Value = "'1.199"
$ws.Range('E39').Value = '  -1.73%  '

# Row 40
$ws.Range('D40').Value = "'8.311"
$ws.Range('E40').Value = '  +3.04%  '

# Row 41
$ws.Range('D41').Value = "'1.010"
$ws.Range('E41').Value = '  +0.34%  '

# Row 42
$ws.Range('D42').Value = "'0.5978"
$ws.Range('E42').Value = '  -0.61%  '

# Row 43
$ws.Range('D43').Value = "'0.1891"
$ws.Range('E43').Value = '  -0.20%  '

# Row 44
$ws.Range('D44').Value = "'10.34"
$ws.Range('E44').Value = '  +0.01%  '

# Row 45
$ws.Range('D45').Value = "'1.250"
$ws.Range('E45').Value = '  -0.85%  '

# Row 46
$ws.Range('D46').Value = "'0.5665"
$ws.Range('E46').Value = '  -1.14%  '

# Row 47
$ws.Range('D47').Value = "'12.17"
$ws.Range('E47').Value = '  -0.51%  '

# Row 48
$ws.Range('D48').Value = "'1.925"
$ws.Range('E48').Value = '  -1.26%  '

# Row 49
$ws.Range('D49').Value = "'0.06764"
$ws.Range('E49').Value = '  -2.23%  '

# Row 50
$ws.Range('D50').Value = "'112.09"
$ws.Range('E50').Value = '  -1.73%  '

# Row 51
$ws.Range('D51').Value = "'3.025"
$ws.Range('E51').Value = '  -11.09%  '
